# Auto-generated edit script updating "想去人数" (column F) counts
# across sheets 展览(1), 演出(2), 本地生活(3), 全部类型(4)
$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = 1; Cell = "F3"; OldValue = 147; NewValue = 150 }
    @{ Sheet = 1; Cell = "F5"; OldValue = 4; NewValue = 5 }
    @{ Sheet = 1; Cell = "F6"; OldValue = 306; NewValue = 308 }
    @{ Sheet = 1; Cell = "F7"; OldValue = 5446; NewValue = 5451 }
    @{ Sheet = 1; Cell = "F9"; OldValue = 7345; NewValue = 7362 }
    @{ Sheet = 1; Cell = "F10"; OldValue = 64; NewValue = 65 }
    @{ Sheet = 1; Cell = "F12"; OldValue = 3753; NewValue = 3758 }
    @{ Sheet = 1; Cell = "F16"; OldValue = 188; NewValue = 189 }
    @{ Sheet = 1; Cell = "F20"; OldValue = 94; NewValue = 97 }
    @{ Sheet = 1; Cell = "F22"; OldValue = 3832; NewValue = 3833 }
    @{ Sheet = 1; Cell = "F24"; OldValue = 5121; NewValue = 5131 }
    @{ Sheet = 1; Cell = "F25"; OldValue = 433; NewValue = 434 }
    @{ Sheet = 1; Cell = "F26"; OldValue = 2043; NewValue = 2046 }
    @{ Sheet = 1; Cell = "F27"; OldValue = 122; NewValue = 125 }
    @{ Sheet = 1; Cell = "F28"; OldValue = 324; NewValue = 326 }
    @{ Sheet = 1; Cell = "F29"; OldValue = 7572; NewValue = 7588 }
    @{ Sheet = 1; Cell = "F30"; OldValue = 28; NewValue = 29 }
    @{ Sheet = 1; Cell = "F34"; OldValue = 147; NewValue = 148 }
    @{ Sheet = 1; Cell = "F35"; OldValue = 1154; NewValue = 1157 }
    @{ Sheet = 1; Cell = "F37"; OldValue = 14; NewValue = 15 }
    @{ Sheet = 1; Cell = "F38"; OldValue = 251; NewValue = 252 }
    @{ Sheet = 1; Cell = "F42"; OldValue = 1169; NewValue = 1170 }
    @{ Sheet = 1; Cell = "F43"; OldValue = 23; NewValue = 24 }
    @{ Sheet = 1; Cell = "F44"; OldValue = 165; NewValue = 167 }
    @{ Sheet = 1; Cell = "F45"; OldValue = 1295; NewValue = 1297 }
    @{ Sheet = 1; Cell = "F46"; OldValue = 1982; NewValue = 1990 }
    @{ Sheet = 1; Cell = "F47"; OldValue = 108; NewValue = 109 }
    @{ Sheet = 1; Cell = "F48"; OldValue = 193; NewValue = 195 }
    @{ Sheet = 2; Cell = "F11"; OldValue = 113; NewValue = 114 }
    @{ Sheet = 3; Cell = "F2"; OldValue = 538; NewValue = 541 }
    @{ Sheet = 3; Cell = "F3"; OldValue = 714; NewValue = 717 }
    @{ Sheet = 4; Cell = "F3"; OldValue = 147; NewValue = 150 }
    @{ Sheet = 4; Cell = "F5"; OldValue = 538; NewValue = 541 }
    @{ Sheet = 4; Cell = "F6"; OldValue = 714; NewValue = 717 }
    @{ Sheet = 4; Cell = "F8"; OldValue = 306; NewValue = 308 }
    @{ Sheet = 4; Cell = "F9"; OldValue = 5446; NewValue = 5451 }
    @{ Sheet = 4; Cell = "F10"; OldValue = 3753; NewValue = 3758 }
    @{ Sheet = 4; Cell = "F14"; OldValue = 188; NewValue = 189 }
    @{ Sheet = 4; Cell = "F17"; OldValue = 94; NewValue = 97 }
    @{ Sheet = 4; Cell = "F21"; OldValue = 3832; NewValue = 3833 }
    @{ Sheet = 4; Cell = "F24"; OldValue = 5121; NewValue = 5131 }
    @{ Sheet = 4; Cell = "F25"; OldValue = 433; NewValue = 434 }
    @{ Sheet = 4; Cell = "F26"; OldValue = 2043; NewValue = 2046 }
    @{ Sheet = 4; Cell = "F27"; OldValue = 122; NewValue = 125 }
    @{ Sheet = 4; Cell = "F28"; OldValue = 324; NewValue = 326 }
    @{ Sheet = 4; Cell = "F29"; OldValue = 7572; NewValue = 7588 }
    @{ Sheet = 4; Cell = "F30"; OldValue = 28; NewValue = 29 }
    @{ Sheet = 4; Cell = "F34"; OldValue = 147; NewValue = 148 }
    @{ Sheet = 4; Cell = "F35"; OldValue = 1154; NewValue = 1157 }
    @{ Sheet = 4; Cell = "F36"; OldValue = 14; NewValue = 15 }
    @{ Sheet = 4; Cell = "F37"; OldValue = 251; NewValue = 252 }
    @{ Sheet = 4; Cell = "F40"; OldValue = 1169; NewValue = 1170 }
    @{ Sheet = 4; Cell = "F41"; OldValue = 23; NewValue = 24 }
    @{ Sheet = 4; Cell = "F42"; OldValue = 165; NewValue = 167 }
    @{ Sheet = 4; Cell = "F44"; OldValue = 1295; NewValue = 1297 }
    @{ Sheet = 4; Cell = "F46"; OldValue = 1982; NewValue = 1990 }
    @{ Sheet = 4; Cell = "F47"; OldValue = 108; NewValue = 109 }
    @{ Sheet = 4; Cell = "F49"; OldValue = 193; NewValue = 195 }
)

$mismatchCount = 0
foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $cell = $ws.Range($edit.Cell)
    if ($cell.Value2 -ne $edit.OldValue) {
        $mismatchCount++
        Write-Output "Warning: sheet $($edit.Sheet) $($edit.Cell) expected $($edit.OldValue) but found $($cell.Value2)"
    }
    $cell.Value = $edit.NewValue
}

Write-Output "Applied $($edits.Count) cell updates ($mismatchCount unexpected prior values)."

$wb.Save()
